# Apply "semana 28 de 2025" update: add column AE (week 28) to the weekly IRA/UCI revision sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell AE1 = "28" (week number), styled like the other header cells (bold + centered).
$ws.Cells.Item(1, 31).Value = 28
$ws.Cells.Item(1, 31).Font.Bold = $true
$ws.Cells.Item(1, 31).HorizontalAlignment = -4108

# Row 28 (CLINICA LOS ROSALES) also had its AC/AD values corrected for this week.
$ws.Cells.Item(28, 29).Value = 1   # AC28
$ws.Cells.Item(28, 30).Value = 0   # AD28 (was 1)

# New weekly counts for column AE (week 28), one value per UPGD row.
$aeValues = @{
    2 = 0
    4 = 0
    5 = 0
    6 = 3
    7 = 0
    8 = 0
    10 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    29 = 0
    30 = 2
    31 = 0
    32 = 0
    34 = 0
    35 = 3
    36 = 0
    37 = 0
    38 = 0
    39 = 0
    40 = 0
    41 = 0
    42 = 0
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 0
    49 = 0
    50 = 0
    52 = 0
    53 = 0
    54 = 0
    55 = 0
    56 = 0
    57 = 0
}
foreach ($row in $aeValues.Keys) {
    $ws.Cells.Item($row, 31).Value = $aeValues[$row]
}

# Row 28 AE value (week 28 count for CLINICA LOS ROSALES).
$ws.Cells.Item(28, 31).Value = 27

Write-Host "Applied week 28 (AE column) update."
